$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3722.25
$ws.Range("I19").Value = 3739.25
$ws.Range("J19").Value = 3713.75
$ws.Range("K19").Value = 3739.25
$ws.Range("L19").Value = 3713.75
$ws.Range("M19").Value = -3564.25
$ws.Range("N19").Value = -4063.75
$ws.Range("H41").Value = 444.375
$ws.Range("I41").Value = 444.375
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 444.375
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4.375
$ws.Range("N41").Value = ""
$ws.Range("H58").Value = 659.2778
$ws.Range("I58").Value = 508.66666
$ws.Range("J58").Value = 809.8889
$ws.Range("K58").Value = 1525.99998
$ws.Range("L58").Value = 2429.6667
$ws.Range("M58").Value = -1375.99998
$ws.Range("N58").Value = -2729.6667
$ws.Range("H62").Value = 5036.9
$ws.Range("I62").Value = 5342.467
$ws.Range("K62").Value = 5342.467
$ws.Range("M62").Value = -4718.467
$ws.Range("H65").Value = 5036.9
$ws.Range("I65").Value = 5342.467
$ws.Range("K65").Value = 26712.335
$ws.Range("M65").Value = -23592.335
$ws.Range("H115").Value = 871
$ws.Range("I115").Value = 871
$ws.Range("K115").Value = 2613
$ws.Range("M115").Value = -1046
$ws.Range("H132").Value = 1044.8
$ws.Range("I132").Value = 1044.8
$ws.Range("K132").Value = 3134.4
$ws.Range("M132").Value = -604.3999999999996
$ws.Range("H135").Value = 943.9487
$ws.Range("I135").Value = 889.8421
$ws.Range("K135").Value = 8008.5789
$ws.Range("M135").Value = -5473.5789
$ws.Range("H137").Value = 5473.304
$ws.Range("J137").Value = 3137.5
$ws.Range("L137").Value = 9412.5
$ws.Range("N137").Value = -14512.5
$ws.Range("H138").Value = 8609.379000000001
$ws.Range("I138").Value = 12416.546
$ws.Range("J138").Value = 6282.778
$ws.Range("K138").Value = 37249.638
$ws.Range("L138").Value = 18848.334
$ws.Range("M138").Value = -32109.638
$ws.Range("N138").Value = -29128.334
$ws.Range("H141").Value = 1645.1034
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21158.23
$ws.Range("I32").Value = 11449.333
$ws.Range("K32").Value = 11449.333
$ws.Range("M32").Value = -11162.333
$ws.Range("H61").Value = 4285
$ws.Range("I61").Value = 4411.8184
$ws.Range("J61").Value = 2890
$ws.Range("K61").Value = 4411.8184
$ws.Range("L61").Value = 2890
$ws.Range("M61").Value = -4199.8184
$ws.Range("N61").Value = -3314
$ws.Range("H110").Value = 3643.682
$ws.Range("I110").Value = 1422.8235
$ws.Range("K110").Value = 1422.8235
$ws.Range("M110").Value = 622.1765
$ws.Range("H132").Value = 1871.1904
$ws.Range("I132").Value = 1796.025
$ws.Range("J132").Value = 3374.5
$ws.Range("K132").Value = 5388.075000000001
$ws.Range("L132").Value = 10123.5
$ws.Range("M132").Value = -2858.075000000001
$ws.Range("N132").Value = -15183.5
$ws.Range("H136").Value = 4285
$ws.Range("I136").Value = 4411.8184
$ws.Range("J136").Value = 2890
$ws.Range("K136").Value = 13235.4552
$ws.Range("L136").Value = 8670
$ws.Range("M136").Value = -10685.4552
$ws.Range("N136").Value = -13770

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5407.857
$ws.Range("I20").Value = 3984.3333
$ws.Range("J20").Value = 7970.2
$ws.Range("K20").Value = 3984.3333
$ws.Range("L20").Value = 7970.2
$ws.Range("M20").Value = -3737.3333
$ws.Range("N20").Value = -8464.200000000001
$ws.Range("H59").Value = 99995
$ws.Range("J59").Value = 99995
$ws.Range("L59").Value = 99995
$ws.Range("N59").Value = -101689
$ws.Range("H105").Value = 2072.6667
$ws.Range("I105").Value = 2407.2
$ws.Range("J105").Value = 400
$ws.Range("K105").Value = 2407.2
$ws.Range("L105").Value = 400
$ws.Range("M105").Value = -660.1999999999998
$ws.Range("N105").Value = -3894
$ws.Range("H134").Value = 5113.0557
$ws.Range("I134").Value = 3156.6924
$ws.Range("K134").Value = 9470.0772
$ws.Range("M134").Value = -6935.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1298.5555
$ws.Range("I6").Value = 651.3077
$ws.Range("K6").Value = 651.3077
$ws.Range("M6").Value = -538.3077
$ws.Range("H12").Value = 1501.6666
$ws.Range("I12").Value = 1501.6666
$ws.Range("K12").Value = 1501.6666
$ws.Range("M12").Value = -1331.6666
$ws.Range("H32").Value = 1100
$ws.Range("I32").Value = 1100
$ws.Range("K32").Value = 1100
$ws.Range("M32").Value = -784
$ws.Range("H58").Value = 6514.5884
$ws.Range("I58").Value = 4273.5
$ws.Range("K58").Value = 4273.5
$ws.Range("M58").Value = -4070.5
$ws.Range("H99").Value = 7126.909
$ws.Range("J99").Value = 7749.625
$ws.Range("L99").Value = 7749.625
$ws.Range("N99").Value = -10745.625
$ws.Range("H122").Value = 3275.7742
$ws.Range("I122").Value = 2770.9583
$ws.Range("K122").Value = 8312.874899999999
$ws.Range("M122").Value = -5862.874899999999
$ws.Range("H126").Value = 7126.909
$ws.Range("J126").Value = 7749.625
$ws.Range("L126").Value = 23248.875
$ws.Range("N126").Value = -28188.875
$ws.Range("H132").Value = 231595.7
$ws.Range("I132").Value = 298361.34
$ws.Range("K132").Value = 895084.02
$ws.Range("M132").Value = -892554.02
$ws.Range("H136").Value = 6514.5884
$ws.Range("I136").Value = 4273.5
$ws.Range("K136").Value = 12820.5
$ws.Range("M136").Value = -10270.5
$ws.Range("H141").Value = 255429.16
$ws.Range("J141").Value = 287870.8
$ws.Range("L141").Value = 287870.8
$ws.Range("N141").Value = -298230.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1267.2069
$ws.Range("J2").Value = 8773.5
$ws.Range("L2").Value = 52641
$ws.Range("N2").Value = -52867
$ws.Range("H22").Value = 1356.6072
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1369.8148
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 4109.4444
$ws.Range("M22").Value = -2831
$ws.Range("N22").Value = -4447.4444
$ws.Range("H27").Value = 1356.6072
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1369.8148
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 4109.4444
$ws.Range("M27").Value = -2898
$ws.Range("N27").Value = -4313.4444
$ws.Range("H117").Value = 1154817.2
$ws.Range("J117").Value = 1667965.1
$ws.Range("L117").Value = 5003895.300000001
$ws.Range("N117").Value = -5010779.300000001
$ws.Range("H122").Value = 29949.875
$ws.Range("J122").Value = 31371.285
$ws.Range("L122").Value = 282341.565
$ws.Range("N122").Value = -287241.565
$ws.Range("H131").Value = 5781
$ws.Range("I131").Value = 3466.6667
$ws.Range("J131").Value = 6475.3
$ws.Range("K131").Value = 10400.0001
$ws.Range("L131").Value = 19425.9
$ws.Range("M131").Value = -5360.000100000001
$ws.Range("N131").Value = -29505.9
$ws.Range("H132").Value = 66669340
$ws.Range("I132").Value = 1596.5
$ws.Range("J132").Value = 76925910
$ws.Range("K132").Value = 14368.5
$ws.Range("L132").Value = 692333190
$ws.Range("M132").Value = -11838.5
$ws.Range("N132").Value = -692338250
$ws.Range("H139").Value = 9330.375
$ws.Range("I139").Value = 3060.5386
$ws.Range("K139").Value = 9181.6158
$ws.Range("M139").Value = -4041.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 19996
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 19996
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 19996
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -20342
$ws.Range("H30").Value = 19996
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 19996
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 19996
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = -20206
$ws.Range("H113").Value = 6345.381
$ws.Range("I113").Value = 6130.231
$ws.Range("K113").Value = 6130.231
$ws.Range("M113").Value = -3960.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 913.7857
$ws.Range("I16").Value = 649.5
$ws.Range("J16").Value = 2499.5
$ws.Range("K16").Value = 649.5
$ws.Range("L16").Value = 2499.5
$ws.Range("M16").Value = -479.5
$ws.Range("N16").Value = -2839.5
$ws.Range("H32").Value = 8146.625
$ws.Range("I32").Value = 8146.625
$ws.Range("K32").Value = 8146.625
$ws.Range("M32").Value = -7829.625
$ws.Range("H122").Value = 8110.875
$ws.Range("I122").Value = 8249.75
$ws.Range("J122").Value = 7972
$ws.Range("K122").Value = 24749.25
$ws.Range("L122").Value = 23916
$ws.Range("M122").Value = -22299.25
$ws.Range("N122").Value = -28816
$ws.Range("H132").Value = 282092
$ws.Range("I132").Value = 303271.38
$ws.Range("K132").Value = 909814.14
$ws.Range("M132").Value = -907284.14
$ws.Range("H136").Value = 275001150
$ws.Range("I136").Value = 171429900
$ws.Range("K136").Value = 514289700
$ws.Range("M136").Value = -514287150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7382.476
$ws.Range("J15").Value = 7393.4287
$ws.Range("L15").Value = 7393.4287
$ws.Range("N15").Value = -7969.4287
$ws.Range("H107").Value = 50662.19
$ws.Range("I107").Value = 2433.1538
$ws.Range("K107").Value = 7299.4614
$ws.Range("M107").Value = -5379.4614
$ws.Range("H132").Value = 307768.66
$ws.Range("I132").Value = 420266
$ws.Range("K132").Value = 1260798
$ws.Range("M132").Value = -1258268
